$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.60286033333333
$ws.Range("H2").Value = 115.808581
$ws.Range("I2").Value = 0.2650212684862838
$ws.Range("J2").Value = 0.2650212684862838
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.787067333333333
$ws.Range("N2").Value = 11.361202
$ws.Range("O2").Value = 0.3088123155238694
$ws.Range("P2").Value = 0.3088123155238694
$ws.Range("Q2").Value = 146.1916313415958
$ws.Range("R2").Value = 1315.724682074362
$ws.Range("S2").Value = 0.08184183158432239
$ws.Range("T2").Value = 0.08184183158432239

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.60286033333333
$ws.Range("H3").Value = 115.808581
$ws.Range("I3").Value = 0.2650212684862838
$ws.Range("J3").Value = 0.2650212684862838
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.323660333333334
$ws.Range("N3").Value = 12.970981
$ws.Range("O3").Value = 0.3525682121685818
$ws.Range("P3").Value = 0.3525682121685817
$ws.Range("Q3").Value = 166.9056559764401
$ws.Range("R3").Value = 1502.150903787961
$ws.Range("S3").Value = 0.09343807481685878
$ws.Range("T3").Value = 0.09343807481685877

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.60286033333333
$ws.Range("H4").Value = 115.808581
$ws.Range("I4").Value = 0.2650212684862838
$ws.Range("J4").Value = 0.2650212684862838
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.152602333333333
$ws.Range("N4").Value = 12.457807
$ws.Range("O4").Value = 0.3386194723075489
$ws.Range("P4").Value = 0.3386194723075489
$ws.Range("Q4").Value = 160.3023278935408
$ws.Range("R4").Value = 1442.720951041867
$ws.Range("S4").Value = 0.08974136208510265
$ws.Range("T4").Value = 0.08974136208510265

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 74.455925
$ws.Range("H5").Value = 223.367775
$ws.Range("I5").Value = 0.5111642898850374
$ws.Range("J5").Value = 0.5111642898850374
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.787067333333333
$ws.Range("N5").Value = 11.361202
$ws.Range("O5").Value = 0.3088123155238694
$ws.Range("P5").Value = 0.3088123155238694
$ws.Range("Q5").Value = 281.9696013406166
$ws.Range("R5").Value = 2537.72641206555
$ws.Range("S5").Value = 0.1578538279725128
$ws.Range("T5").Value = 0.1578538279725128

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 74.455925
$ws.Range("H6").Value = 223.367775
$ws.Range("I6").Value = 0.5111642898850374
$ws.Range("J6").Value = 0.5111642898850374
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.323660333333334
$ws.Range("N6").Value = 12.970981
$ws.Range("O6").Value = 0.3525682121685818
$ws.Range("P6").Value = 0.3525682121685817
$ws.Range("Q6").Value = 321.9221295041417
$ws.Range("R6").Value = 2897.299165537275
$ws.Range("S6").Value = 0.1802202798091903
$ws.Range("T6").Value = 0.1802202798091903

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 74.455925
$ws.Range("H7").Value = 223.367775
$ws.Range("I7").Value = 0.5111642898850374
$ws.Range("J7").Value = 0.5111642898850374
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.152602333333333
$ws.Range("N7").Value = 12.457807
$ws.Range("O7").Value = 0.3386194723075489
$ws.Range("P7").Value = 0.3386194723075489
$ws.Range("Q7").Value = 309.1858478854916
$ws.Range("R7").Value = 2782.672630969425
$ws.Range("S7").Value = 0.1730901821033343
$ws.Range("T7").Value = 0.1730901821033343

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 32.60069533333333
$ws.Range("H8").Value = 97.80208599999999
$ws.Range("I8").Value = 0.2238144416286788
$ws.Range("J8").Value = 0.2238144416286788
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.787067333333333
$ws.Range("N8").Value = 11.361202
$ws.Range("O8").Value = 0.3088123155238694
$ws.Range("P8").Value = 0.3088123155238694
$ws.Range("Q8").Value = 123.4610283408191
$ws.Range("R8").Value = 1111.149255067372
$ws.Range("S8").Value = 0.06911665596703419
$ws.Range("T8").Value = 0.0691166559670342

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 32.60069533333333
$ws.Range("H9").Value = 97.80208599999999
$ws.Range("I9").Value = 0.2238144416286788
$ws.Range("J9").Value = 0.2238144416286788
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.323660333333334
$ws.Range("N9").Value = 12.970981
$ws.Range("O9").Value = 0.3525682121685818
$ws.Range("P9").Value = 0.3525682121685817
$ws.Range("Q9").Value = 140.9543332518184
$ws.Range("R9").Value = 1268.588999266366
$ws.Range("S9").Value = 0.07890985754253267
$ws.Range("T9").Value = 0.07890985754253266

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 32.60069533333333
$ws.Range("H10").Value = 97.80208599999999
$ws.Range("I10").Value = 0.2238144416286788
$ws.Range("J10").Value = 0.2238144416286788
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.152602333333333
$ws.Range("N10").Value = 12.457807
$ws.Range("O10").Value = 0.3386194723075489
$ws.Range("P10").Value = 0.3386194723075489
$ws.Range("Q10").Value = 135.3777235094891
$ws.Range("R10").Value = 1218.399511585402
$ws.Range("S10").Value = 0.07578792811911189
$ws.Range("T10").Value = 0.0757879281191119
